# Applies the "changement mineurs sur le xls" edit:
#   - C13: "Le logiciel est moins plus flexible pour l'ecriture des mains."
#          -> "Le logiciel est plus flexible pour l'ecriture des mains."
#   - E3:  "...le logiciel gere les erreurs de saisie et affiche..."
#          -> "...le logiciel gere le cas ou une meme carte est utilisee deux fois, et affiche..."
#   - Selection moves from E14 to G2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C13").Value = "Le logiciel est plus flexible pour l'écriture des mains."

$ws.Range("E3").Value = "Le logiciel reconnait deux mains de cinq cartes entrées par l'utilisateur et indique laquelle est la plus forte sans tenir compte des combinaisons: couleur et suite. Le logiciel ne gère qu'un paquet contenant que la famille pique, le logiciel gère le cas où une même carte est utilisée deux fois, et affiche le résultat comme les spécifications.  La saisie de carte est plus simple, les têtes peuvent être tapées avec leur initiale."

$ws.Range("G2").Select()
